$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the existing "总计" sheet,
#    and keep a freshly-built "总计" sheet after it.
#
#    To keep the sheet-level formatting ("总计"'s sheetPr/pageMargins/
#    sheetFormatPr) on both resulting sheets, we duplicate the current
#    "总计" sheet (so the copy inherits all of its page/format settings),
#    then rename the original to "2022-Q1" and the copy back to "总计".
#    This also reproduces the target sheetId numbering: the original sheet
#    (sheetId=6) becomes "2022-Q1", and the new copy (sheetId=7) becomes
#    "总计".
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($null, $total)
$totalCopy = $wb.Worksheets.Item($total.Index + 1)

$total.Name = "2022-Q1"
$totalCopy.Name = "总计"

$q1 = $total

# ---------------------------------------------------------------------------
# 2) Rebuild the "2022-Q1" sheet as a fund-holdings detail sheet (same shape
#    as the other quarterly sheets, e.g. "2021-Q4"). Reuse that sheet's
#    header/index-column formatting via copy/paste-special so styles match
#    the rest of the workbook.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")

# Wipe the old "总计" data/format still sitting in this sheet's cells, then
# lay down the 7-column fund-holdings header formatting from the template.
$q1.Cells.Clear()

$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "000893"
$q1.Range("C2").Value = "工银创新动力股票"
$q1.Range("D2:G2").NumberFormat = "@"
$q1.Range("D2").Value = "11.07"
$q1.Range("E2").Value = "81.01"
$q1.Range("F2").Value = "4.02"
$q1.Range("G2").Value = "0.4450"
$q1.Range("H2").Value = 4

$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "160639"
$q1.Range("C3").Value = "鹏华中证高铁产业指数（LOF）"
$q1.Range("D3:G3").NumberFormat = "@"
$q1.Range("D3").Value = "0.89"
$q1.Range("E3").Value = "94.72"
$q1.Range("F3").Value = "2.28"
$q1.Range("G3").Value = "0.0203"
$q1.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 3) Rebuild the "总计" sheet: same 3-column rollup table as before, with a
#    new "2022-Q1" row inserted at the top and every other row shifted down.
# ---------------------------------------------------------------------------
$totalCopy.Cells.Clear()

$template.Range("B1:D1").Copy()
$totalCopy.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$totalCopy.Range("A2:A7").PasteSpecial(-4122)

$totalCopy.Range("B1").Value = "日期"
$totalCopy.Range("C1").Value = "持有数量(只)"
$totalCopy.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @(0, "2022-Q1", 2, 0.47),
    @(1, "2021-Q4", 1, 0.32),
    @(2, "2021-Q3", 9, 0.89),
    @(3, "2021-Q2", 5, 0.5600000000000001),
    @(4, "2021-Q1", 4, 0.42),
    @(5, "2020-Q4", 5, 0.53)
)

$r = 2
foreach ($row in $rows) {
    $totalCopy.Range("A$r").Value = $row[0]
    $totalCopy.Range("B$r").Value = $row[1]
    $totalCopy.Range("C$r").Value = $row[2]
    $totalCopy.Range("D$r").Value = $row[3]
    $r++
}
